$wb = $excel.ActiveWorkbook

# A new handoff was generated for the "5fdc99be-..." file, refreshing its
# "Latest Handoff Date" / "Latest Handoff Datetime" timestamps in the
# Overview sheet and in each per-locale detail sheet.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-23 12:43:28"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-23 12:43:24"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-23 12:43:28"
